$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template style source cells (existing row 50 has the established formatting for this table)

# ===== Row 51 =====
$ws.Range("A50").Copy()
$ws.Range("A51").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B51:E51").PasteSpecial(-4122)
$ws.Range("BA51:BG51").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F51").PasteSpecial(-4122)

$ws.Range("A51").Value = 44571.71661831018
$ws.Range("B51").Value = "Enrymar Cisneros"
$ws.Range("C51").Value = "<40"
$ws.Range("D51").Value = "2013CISN01"
$ws.Range("E51").Value = "Megaminx"
$ws.Range("F51").Value = "https://www.facebook.com/events/343359980546742/?post_id=350246939858046&view=permalink"
$ws.Range("BA51").Value = "1:34.36"
$ws.Range("BB51").Value = "1:42.12"
$ws.Range("BC51").Value = "1:39.64"
$ws.Range("BD51").Value = "1:38.92"
$ws.Range("BE51").Value = "1:36.72"
$ws.Range("BF51").Value = "1:34.36"
$ws.Range("BG51").Value = "1:38.43"

# ===== Row 52 =====
$ws.Range("A50").Copy()
$ws.Range("A52").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B52:E52").PasteSpecial(-4122)
$ws.Range("H52:N52").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F52").PasteSpecial(-4122)

$ws.Range("A52").Value = 44571.71736016204
$ws.Range("B52").Value = "Enrymar Cisneros"
$ws.Range("C52").Value = "<40"
$ws.Range("D52").Value = "2013CISN01"
$ws.Range("E52").Value = "2x2x2"
$ws.Range("F52").Value = "https://www.facebook.com/events/343359980546742/?post_id=350165623199511&view=permalink"
$ws.Range("H52").Value = 4.74
$ws.Range("I52").Value = 8.17
$ws.Range("J52").Value = 5.35
$ws.Range("K52").Value = 5.87
$ws.Range("L52").Value = 7.6
$ws.Range("M52").Value = 4.74
$ws.Range("N52").Value = 6.27

# ===== Row 53 =====
$ws.Range("A50").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B53:E53").PasteSpecial(-4122)
$ws.Range("O53:U53").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F53").PasteSpecial(-4122)

$ws.Range("A53").Value = 44571.71794107639
$ws.Range("B53").Value = "Enrymar Cisneros"
$ws.Range("C53").Value = "<40"
$ws.Range("D53").Value = "2013CISN01"
$ws.Range("E53").Value = "3x3x3"
$ws.Range("F53").Value = "https://www.facebook.com/events/343359980546742/?post_id=350163353199738&view=permalink"
$ws.Range("O53").Value = 12.35
$ws.Range("P53").Value = 13.93
$ws.Range("Q53").Value = 12.61
$ws.Range("R53").Value = 12.2
$ws.Range("S53").Value = 12.2
$ws.Range("T53").Value = 12.2
$ws.Range("U53").Value = 12.39

# ===== Row 54 =====
$ws.Range("A50").Copy()
$ws.Range("A54").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B54:E54").PasteSpecial(-4122)
$ws.Range("BH54:BN54").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F54").PasteSpecial(-4122)

$ws.Range("A54").Value = 44571.719053622684
$ws.Range("B54").Value = "Enrymar Cisneros"
$ws.Range("C54").Value = "<40"
$ws.Range("D54").Value = "2013CISN01"
$ws.Range("E54").Value = "Pyraminx"
$ws.Range("F54").Value = "https://www.facebook.com/events/1083505512394794/?post_id=1091527064925972&view=permalink"
$ws.Range("BH54").Value = 9.2
$ws.Range("BI54").Value = 7.6
$ws.Range("BJ54").Value = 7.1
$ws.Range("BK54").Value = 9.53
$ws.Range("BL54").Value = 5.99
$ws.Range("BM54").Value = 5.99
$ws.Range("BN54").Value = 7.97

# ===== Row 55 =====
$ws.Range("A50").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B55:E55").PasteSpecial(-4122)
$ws.Range("BO55:BU55").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F55").PasteSpecial(-4122)

$ws.Range("A55").Value = 44571.72001913194
$ws.Range("B55").Value = "Enrymar Cisneros"
$ws.Range("C55").Value = "<40"
$ws.Range("D55").Value = "2013CISN01"
$ws.Range("E55").Value = "Skewb"
$ws.Range("F55").Value = "https://www.facebook.com/events/1083505512394794/?post_id=1091524834926195&view=permalink"
$ws.Range("BO55").Value = 3.85
$ws.Range("BP55").Value = 13.41
$ws.Range("BQ55").Value = 6.75
$ws.Range("BR55").Value = 8.49
$ws.Range("BS55").Value = 8.61
$ws.Range("BT55").Value = 3.85
$ws.Range("BU55").Value = 7.95

# ===== Row 56 =====
$ws.Range("A50").Copy()
$ws.Range("A56").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B56:E56").PasteSpecial(-4122)
$ws.Range("AO56:AS56").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F56").PasteSpecial(-4122)

$ws.Range("A56").Value = 44571.72069488426
$ws.Range("B56").Value = "Enrymar Cisneros"
$ws.Range("C56").Value = "<40"
$ws.Range("D56").Value = "2013CISN01"
$ws.Range("E56").Value = "7x7x7"
$ws.Range("F56").Value = "https://www.facebook.com/events/364077578855426/?post_id=371890168074167&view=permalink"
$ws.Range("AO56").Value = "5:45.96"
$ws.Range("AP56").Value = "5:25.99"
$ws.Range("AQ56").Value = "5:40.08"
$ws.Range("AR56").Value = "5:25.99"
$ws.Range("AS56").Value = "5:37.34"

# ===== Row 57 =====
$ws.Range("A50").Copy()
$ws.Range("A57").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B57:E57").PasteSpecial(-4122)
$ws.Range("AJ57:AN57").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F57").PasteSpecial(-4122)

$ws.Range("A57").Value = 44571.72126261574
$ws.Range("B57").Value = "Enrymar Cisneros"
$ws.Range("C57").Value = "<40"
$ws.Range("D57").Value = "2013CISN01"
$ws.Range("E57").Value = "6x6x6"
$ws.Range("F57").Value = "https://www.facebook.com/events/364077578855426/?post_id=371704864759364&view=permalink"
$ws.Range("AJ57").Value = "3:40.63"
$ws.Range("AK57").Value = "3:39.47"
$ws.Range("AL57").Value = "3:22.95"
$ws.Range("AM57").Value = "3:22.95"
$ws.Range("AN57").Value = "3:34.35"

# ===== Row 58 =====
$ws.Range("A50").Copy()
$ws.Range("A58").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B58:E58").PasteSpecial(-4122)
$ws.Range("V58:AB58").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F58").PasteSpecial(-4122)

$ws.Range("A58").Value = 44571.7220156713
$ws.Range("B58").Value = "Enrymar Cisneros"
$ws.Range("C58").Value = "<40"
$ws.Range("D58").Value = "2013CISN01"
$ws.Range("E58").Value = "4x4x4"
$ws.Range("F58").Value = "https://www.facebook.com/events/364077578855426/?post_id=371693138093870&view=permalink"
$ws.Range("V58").Value = 53.59
$ws.Range("W58").Value = 51.09
$ws.Range("X58").Value = 52.06
$ws.Range("Y58").Value = 54.0
$ws.Range("Z58").Value = "1:06.77"
$ws.Range("AA58").Value = 51.09
$ws.Range("AB58").Value = 53.22

# ===== Row 59 =====
$ws.Range("A50").Copy()
$ws.Range("A59").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("B59:E59").PasteSpecial(-4122)
$ws.Range("AC59:AI59").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F59").PasteSpecial(-4122)

$ws.Range("A59").Value = 44571.722717743054
$ws.Range("B59").Value = "Enrymar Cisneros"
$ws.Range("C59").Value = "<40"
$ws.Range("D59").Value = "2013CISN01"
$ws.Range("E59").Value = "5x5x5"
$ws.Range("F59").Value = "https://www.facebook.com/events/364077578855426/?post_id=371689508094233&view=permalink"
$ws.Range("AC59").Value = "1:56.90"
$ws.Range("AD59").Value = "1:45.00"
$ws.Range("AE59").Value = "1:46.67"
$ws.Range("AF59").Value = "1:52.90"
$ws.Range("AG59").Value = "1:47.39"
$ws.Range("AH59").Value = "1:45.00"
$ws.Range("AI59").Value = "1:48.99"

# ===== Hyperlinks for the Facebook link column (F) =====
$ws.Hyperlinks.Add($ws.Range("F51"), "https://www.facebook.com/events/343359980546742/?post_id=350246939858046&view=permalink")
$ws.Hyperlinks.Add($ws.Range("F52"), "https://www.facebook.com/events/343359980546742/?post_id=350165623199511&view=permalink")
$ws.Hyperlinks.Add($ws.Range("F53"), "https://www.facebook.com/events/343359980546742/?post_id=350163353199738&view=permalink")
$ws.Hyperlinks.Add($ws.Range("F54"), "https://www.facebook.com/events/1083505512394794/?post_id=1091527064925972&view=permalink")
$ws.Hyperlinks.Add($ws.Range("F55"), "https://www.facebook.com/events/1083505512394794/?post_id=1091524834926195&view=permalink")
$ws.Hyperlinks.Add($ws.Range("F56"), "https://www.facebook.com/events/364077578855426/?post_id=371890168074167&view=permalink")
$ws.Hyperlinks.Add($ws.Range("F57"), "https://www.facebook.com/events/364077578855426/?post_id=371704864759364&view=permalink")
$ws.Hyperlinks.Add($ws.Range("F58"), "https://www.facebook.com/events/364077578855426/?post_id=371693138093870&view=permalink")
$ws.Hyperlinks.Add($ws.Range("F59"), "https://www.facebook.com/events/364077578855426/?post_id=371689508094233&view=permalink")

# Re-apply the established link-cell formatting (Hyperlinks.Add swaps in the built-in "Hyperlink" style)
$ws.Range("F50").Copy()
$ws.Range("F51").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F52").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F53").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F54").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F55").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F56").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F57").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F58").PasteSpecial(-4122)
$ws.Range("F50").Copy()
$ws.Range("F59").PasteSpecial(-4122)
